$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header values in P1 and Q1, matching the style used by B1:O1 (bold, bordered, centered)
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Update existing data columns I, K, M, O for rows 2-25
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # column I
    $ws.Cells.Item($r, 11).Value = 1  # column K
    $ws.Cells.Item($r, 13).Value = 2  # column M
    $ws.Cells.Item($r, 15).Value = 1  # column O

    # New columns P and Q with value 2, no special style
    $ws.Cells.Item($r, 16).Value = 2  # column P
    $ws.Cells.Item($r, 17).Value = 2  # column Q
}
